# Scheduled runner update: refresh Tonberry_Profits leve profit/price figures
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets (market-price driven columns H-N).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2137.25
$ws.Range("I62").Value = 1274.5
$ws.Range("K62").Value = 1274.5
$ws.Range("M62").Value = -650.5
$ws.Range("H65").Value = 2137.25
$ws.Range("I65").Value = 1274.5
$ws.Range("K65").Value = 6372.5
$ws.Range("M65").Value = -3252.5
$ws.Range("H112").Value = 7062.5
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 7062.5
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 21187.5
$ws.Range("M112").ClearContents()
$ws.Range("N112").Value = -23403.5
$ws.Range("H132").Value = 1196.3243
$ws.Range("I132").Value = 998.53125
$ws.Range("J132").Value = 2462.2
$ws.Range("K132").Value = 2995.59375
$ws.Range("L132").Value = 7386.599999999999
$ws.Range("M132").Value = -465.59375
$ws.Range("N132").Value = -12446.6
$ws.Range("H137").Value = 1945
$ws.Range("I137").Value = 1707.1428
$ws.Range("J137").Value = 2500
$ws.Range("K137").Value = 5121.428400000001
$ws.Range("L137").Value = 7500
$ws.Range("M137").Value = -2571.428400000001
$ws.Range("N137").Value = -12600
$ws.Range("H138").Value = 2406
$ws.Range("J138").Value = 2349.75
$ws.Range("L138").Value = 7049.25
$ws.Range("N138").Value = -17329.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2930.9155
$ws.Range("I32").Value = 2018.1052
$ws.Range("J32").Value = 6647.357
$ws.Range("K32").Value = 2018.1052
$ws.Range("L32").Value = 6647.357
$ws.Range("M32").Value = -1731.1052
$ws.Range("N32").Value = -7221.357
$ws.Range("H61").Value = 4184.7
$ws.Range("I61").Value = 2135.1667
$ws.Range("K61").Value = 2135.1667
$ws.Range("M61").Value = -1923.1667
$ws.Range("H63").Value = 8149.75
$ws.Range("I63").Value = 8149.75
$ws.Range("K63").Value = 8149.75
$ws.Range("M63").Value = -7463.75
$ws.Range("H66").Value = 8149.75
$ws.Range("I66").Value = 8149.75
$ws.Range("K66").Value = 40748.75
$ws.Range("M66").Value = -37316.75
$ws.Range("H132").Value = 1593.4509
$ws.Range("I132").Value = 1236.0555
$ws.Range("J132").Value = 2451.2
$ws.Range("K132").Value = 3708.1665
$ws.Range("L132").Value = 7353.599999999999
$ws.Range("M132").Value = -1178.1665
$ws.Range("N132").Value = -12413.6
$ws.Range("H136").Value = 4184.7
$ws.Range("I136").Value = 2135.1667
$ws.Range("K136").Value = 6405.500100000001
$ws.Range("M136").Value = -3855.500100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 1765
$ws.Range("I7").Value = 1765
$ws.Range("K7").Value = 1765
$ws.Range("M7").Value = -1652
$ws.Range("H20").Value = 1981.2727
$ws.Range("I20").Value = 2029.6
$ws.Range("J20").Value = 1498
$ws.Range("K20").Value = 2029.6
$ws.Range("L20").Value = 1498
$ws.Range("M20").Value = -1782.6
$ws.Range("N20").Value = -1992
$ws.Range("H134").Value = 21795
$ws.Range("I134").Value = 21795
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 65385
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -62850
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2172.8235
$ws.Range("I31").Value = 1796.2
$ws.Range("K31").Value = 1796.2
$ws.Range("M31").Value = -1501.2
$ws.Range("H34").Value = 2172.8235
$ws.Range("I34").Value = 1796.2
$ws.Range("K34").Value = 1796.2
$ws.Range("M34").Value = -1594.2
$ws.Range("H134").Value = 1718.303
$ws.Range("I134").Value = 1543.25
$ws.Range("J134").Value = 2698.6
$ws.Range("K134").Value = 4629.75
$ws.Range("L134").Value = 8095.799999999999
$ws.Range("M134").Value = -2094.75
$ws.Range("N134").Value = -13165.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 7143014
$ws.Range("I4").Value = 174.75
$ws.Range("J4").Value = 16666800
$ws.Range("K4").Value = 524.25
$ws.Range("L4").Value = 50000400
$ws.Range("M4").Value = -412.25
$ws.Range("N4").Value = -50000624
$ws.Range("H5").Value = 797.5454999999999
$ws.Range("J5").Value = 860.375
$ws.Range("L5").Value = 2581.125
$ws.Range("N5").Value = -2805.125
$ws.Range("H59").Value = 1000
$ws.Range("J59").Value = 1000
$ws.Range("L59").Value = 3000
$ws.Range("N59").Value = -4080
$ws.Range("H104").Value = 3042.7827
$ws.Range("J104").Value = 3436.5789
$ws.Range("L104").Value = 10309.7367
$ws.Range("N104").Value = -15551.7367
$ws.Range("H122").Value = 1069.9286
$ws.Range("J122").Value = 1109.091
$ws.Range("L122").Value = 9981.819
$ws.Range("N122").Value = -14881.819
$ws.Range("H131").Value = 783.61
$ws.Range("J131").Value = 800.3617
$ws.Range("L131").Value = 2401.0851
$ws.Range("N131").Value = -12481.0851
$ws.Range("H135").Value = 797.5454999999999
$ws.Range("J135").Value = 860.375
$ws.Range("L135").Value = 7743.375
$ws.Range("N135").Value = -12813.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2250
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 2250
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 2250
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -4246
$ws.Range("H83").Value = 2250
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 2250
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 11250
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -21234
$ws.Range("H122").Value = 1936.6111
$ws.Range("I122").Value = 1487.3636
$ws.Range("K122").Value = 4462.0908
$ws.Range("M122").Value = -2012.0908
$ws.Range("H132").Value = 1926813.2
$ws.Range("I132").Value = 2962250.2
$ws.Range("J132").Value = 3858.8572
$ws.Range("K132").Value = 8886750.600000001
$ws.Range("L132").Value = 11576.5716
$ws.Range("M132").Value = -8884220.600000001
$ws.Range("N132").Value = -16636.5716

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3152.4666
$ws.Range("J7").Value = 4160.8335
$ws.Range("L7").Value = 4160.8335
$ws.Range("N7").Value = -4384.8335
$ws.Range("H126").Value = 3152.4666
$ws.Range("J126").Value = 4160.8335
$ws.Range("L126").Value = 12482.5005
$ws.Range("N126").Value = -17422.5005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 1847
$ws.Range("J14").Value = 1847
$ws.Range("L14").Value = 1847
$ws.Range("N14").Value = -2183
$ws.Range("H136").Value = 55558716
$ws.Range("I136").Value = 138891360
$ws.Range("J136").Value = 3616.3333
$ws.Range("K136").Value = 416674080
$ws.Range("L136").Value = 10848.9999
$ws.Range("M136").Value = -416671530
$ws.Range("N136").Value = -15948.9999
